$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "26.871.42"
Set-TextCell $ws.Range("E2") "  -2.29%  "

# Row 3
Set-TextCell $ws.Range("D3") "1.830.93"
Set-TextCell $ws.Range("E3") "  -2.08%  "

# Row 4
Set-TextCell $ws.Range("D4") "1.007"
Set-TextCell $ws.Range("E4") "  +0.28%  "

# Row 5
Set-TextCell $ws.Range("D5") "310.27"
Set-TextCell $ws.Range("E5") "  -1.92%  "

# Row 6
Set-TextCell $ws.Range("E6") "  +0.18%  "

# Row 7
Set-TextCell $ws.Range("D7") "0.4615"
Set-TextCell $ws.Range("E7") "  -1.22%  "

# Row 8
Set-TextCell $ws.Range("D8") "0.3652"
Set-TextCell $ws.Range("E8") "  -2.46%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.07166"
Set-TextCell $ws.Range("E9") "  -3.21%  "

# Row 10
Set-TextCell $ws.Range("D10") "0.8790"
Set-TextCell $ws.Range("E10") "  -1.27%  "

# Row 11
Set-TextCell $ws.Range("D11") "0.07832"
Set-TextCell $ws.Range("E11") "  -1.66%  "

# Row 12
Set-TextCell $ws.Range("D12") "19.57"
Set-TextCell $ws.Range("E12") "  -2.59%  "

# Row 13
Set-TextCell $ws.Range("D13") "1.875.84"
Set-TextCell $ws.Range("E13") "  +1.51%  "

# Row 14
Set-TextCell $ws.Range("D14") "5.324"
Set-TextCell $ws.Range("E14") "  -2.18%  "

# Row 15
Set-TextCell $ws.Range("D15") "6.386"
Set-TextCell $ws.Range("E15") "  -3.42%  "

# Row 16
Set-TextCell $ws.Range("D16") "88.56"
Set-TextCell $ws.Range("E16") "  -4.69%  "

# Row 17
Set-TextCell $ws.Range("E17") "  +0.14%  "

# Row 18
Set-TextCell $ws.Range("D18") "0.000008737"
Set-TextCell $ws.Range("E18") "  -2.65%  "

# Row 19
Set-TextCell $ws.Range("D19") "1.006"
Set-TextCell $ws.Range("E19") "  +0.19%  "

# Row 20
Set-TextCell $ws.Range("D20") "26.897.76"
Set-TextCell $ws.Range("E20") "  -2.28%  "

# Row 21
Set-TextCell $ws.Range("D21") "14.50"
Set-TextCell $ws.Range("E21") "  -3.21%  "

# Row 22
Set-TextCell $ws.Range("D22") "4.997"
Set-TextCell $ws.Range("E22") "  -3.84%  "

# Row 23
Set-TextCell $ws.Range("E23") "  -1.73%  "

# Row 24
Set-TextCell $ws.Range("D24") "1.977"
Set-TextCell $ws.Range("E24") "  +4.87%  "

# Row 25
Set-TextCell $ws.Range("D25") "150.79"
Set-TextCell $ws.Range("E25") "  -1.54%  "

# Row 26
Set-TextCell $ws.Range("D26") "18.21"
Set-TextCell $ws.Range("E26") "  -2.18%  "

# Row 27
Set-TextCell $ws.Range("D27") "1.994"
Set-TextCell $ws.Range("E27") "  -5.04%  "

# Row 28
Set-TextCell $ws.Range("D28") "113.52"
Set-TextCell $ws.Range("E28") "  -3.47%  "

# Row 29
Set-TextCell $ws.Range("D29") "4.943"
Set-TextCell $ws.Range("E29") "  -4.63%  "

# Row 30
Set-TextCell $ws.Range("D30") "0.08816"
Set-TextCell $ws.Range("E30") "  -1.15%  "

# Row 31
Set-TextCell $ws.Range("D31") "3.097"
Set-TextCell $ws.Range("E31") "  +2.69%  "

# Row 32
Set-TextCell $ws.Range("D32") "0.7626"
Set-TextCell $ws.Range("E32") "  +0.80%  "

# Row 33
Set-TextCell $ws.Range("D33") "4.462"
Set-TextCell $ws.Range("E33") "  -0.95%  "

# Row 34
Set-TextCell $ws.Range("E34") "  -2.31%  "

# Row 35
Set-TextCell $ws.Range("D35") "2.655"
Set-TextCell $ws.Range("E35") "  -0.78%  "

# Row 36
Set-TextCell $ws.Range("D36") "1.091"
Set-TextCell $ws.Range("E36") "  +0.50%  "

# Row 37
Set-TextCell $ws.Range("D37") "0.01925"
Set-TextCell $ws.Range("E37") "  -2.43%  "

# Row 38
Set-TextCell $ws.Range("D38") "0.05136"
Set-TextCell $ws.Range("E38") "  -3.36%  "

# Row 39
Set-TextCell $ws.Range("D39") "2.920"
Set-TextCell $ws.Range("E39") "  -2.42%  "

# Row 40
Set-TextCell $ws.Range("D40") "6.948"
Set-TextCell $ws.Range("E40") "  -3.73%  "

# Row 41
Set-TextCell $ws.Range("D41") "0.4982"
Set-TextCell $ws.Range("E41") "  -5.22%  "

# Row 42
Set-TextCell $ws.Range("D42") "0.1595"
Set-TextCell $ws.Range("E42") "  -3.36%  "

# Row 43
Set-TextCell $ws.Range("D43") "8.365"

# Row 44
Set-TextCell $ws.Range("B44") "EnergySwap"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D44") "10.24"
Set-TextCell $ws.Range("E44") "  -1.01%  "

# Row 45
Set-TextCell $ws.Range("B45") "Decentraland"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell $ws.Range("D45") "0.4666"
Set-TextCell $ws.Range("E45") "  -5.17%  "

# Row 46
Set-TextCell $ws.Range("B46") "PaxDollar"
Set-TextCell $ws.Range("C46") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell $ws.Range("D46") "1.006"
Set-TextCell $ws.Range("E46") "  +0.17%  "

# Row 47
Set-TextCell $ws.Range("D47") "102.50"
Set-TextCell $ws.Range("E47") "  -1.52%  "

# Row 48
Set-TextCell $ws.Range("D48") "1.609"
Set-TextCell $ws.Range("E48") "  -4.08%  "

# Row 49
Set-TextCell $ws.Range("D49") "0.06100"
Set-TextCell $ws.Range("E49") "  -2.59%  "

# Row 50
Set-TextCell $ws.Range("D50") "64.61"
Set-TextCell $ws.Range("E50") "  -2.14%  "

# Row 51
Set-TextCell $ws.Range("D51") "36.31"
Set-TextCell $ws.Range("E51") "  -2.77%  "
